# Repoint the EoDSDwSP formulas away from the (about to be removed) "Texas
# Notes" sheet and back to the original "Calculations" sheet, then drop the
# "Texas Notes" worksheet entirely.
$wb = $excel.ActiveWorkbook

$eods = $wb.Worksheets.Item("EoDSDwSP")
$eods.Range("B2").Formula = "=Calculations!B9"
$eods.Range("B4").Formula = "=Calculations!B10"

$texasNotes = $wb.Worksheets.Item("Texas Notes")
$texasNotes.Delete() | Out-Null

# Remove the hyperlink that was added to the About sheet (cell B6) along
# with the "Hyperlink" cell style it pulled in.
$about = $wb.Worksheets.Item("About")
$about.Range("B6").Hyperlinks.Delete()
$about.Range("B6").Style = "Normal"
$wb.Styles.Item("Hyperlink").Delete()

# Restore the prior view state: EoDSDwSP's selection back to B2 and no
# longer the active tab, with About active/selected at A12 instead.
# (Re-fetch the sheet objects: the worksheet collection changed above, so
# the earlier $eods / $about references are stale.)
$eods = $wb.Worksheets.Item("EoDSDwSP")
$eods.Activate() | Out-Null
$eods.Range("B2").Select() | Out-Null

$about = $wb.Worksheets.Item("About")
$about.Activate() | Out-Null
$about.Range("A12").Select() | Out-Null
